$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "competition a" -> "EG"
$ws.Name = "EG"

# Update the header info block
$ws.Range("B1").Value = "EG"
$ws.Range("B2").Value = "https://noon.com"
$ws.Range("B3").Value = 44553
$ws.Range("B4").Value = 4

# Remove the two data rows (7 and 8) that listed the students
$ws.Rows("7:8").Delete()

# Narrow column B slightly (stored width target 17.0078125 chars;
# the COM width model here quantizes to 1/6-character steps, so feed it
# the value whose quantized result lands closest/exactly on that target)
$ws.Columns("B").ColumnWidth = 17.0078125 - 5/6
